# Applies the data refresh to the "展览" (sheet 1) and "全部类型" (sheet 4)
# worksheets, which both hold the same table of 漫展/演出 events.
#
# For each affected row:
#   - column F ("想去人数") gets an updated interest count
#   - column G ("最低票价") gets an updated minimum price, and for two rows
#     the numeric price becomes the text "不可售" (not currently for sale)

function Apply-Updates($worksheet, $updates) {
    foreach ($row in $updates.Keys) {
        $vals = $updates[$row]
        $worksheet.Cells.Item($row, 6).Value = $vals.F
        $worksheet.Cells.Item($row, 7).Value = $vals.G
    }
}

$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet index 1)
$sheet1Updates = @{
    2  = @{ F = 271;   G = 20 }
    3  = @{ F = 587;   G = 49.6 }
    4  = @{ F = 6923;  G = 80 }
    5  = @{ F = 102;   G = 60 }
    7  = @{ F = 182;   G = "不可售" }
    8  = @{ F = 68;    G = "不可售" }
    9  = @{ F = 1139;  G = 58 }
    10 = @{ F = 16479; G = 60 }
    14 = @{ F = 352;   G = 50 }
    15 = @{ F = 196;   G = 20 }
    17 = @{ F = 11492; G = 60 }
    18 = @{ F = 20;    G = 39.9 }
    19 = @{ F = 1138;  G = 60 }
    20 = @{ F = 4537;  G = 39 }
    21 = @{ F = 384;   G = 70 }
    22 = @{ F = 393;   G = 75 }
    24 = @{ F = 863;   G = 60 }
    26 = @{ F = 145;   G = 58 }
}

# "全部类型" sheet (sheet index 4) — same events, slightly different row numbers
$sheet4Updates = @{
    2  = @{ F = 271;   G = 20 }
    3  = @{ F = 587;   G = 49.6 }
    4  = @{ F = 6923;  G = 80 }
    5  = @{ F = 102;   G = 60 }
    7  = @{ F = 182;   G = "不可售" }
    8  = @{ F = 68;    G = "不可售" }
    10 = @{ F = 1139;  G = 58 }
    11 = @{ F = 16479; G = 60 }
    15 = @{ F = 352;   G = 50 }
    16 = @{ F = 196;   G = 20 }
    20 = @{ F = 11492; G = 60 }
    21 = @{ F = 20;    G = 39.9 }
    22 = @{ F = 1138;  G = 60 }
    23 = @{ F = 4537;  G = 39 }
    24 = @{ F = 384;   G = 70 }
    25 = @{ F = 393;   G = 75 }
    27 = @{ F = 863;   G = 60 }
    29 = @{ F = 145;   G = 58 }
}

$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

Apply-Updates $ws1 $sheet1Updates
Apply-Updates $ws4 $sheet4Updates
